$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.110.39"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "'1.850.92"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.62%  "
$ws.Range("D5").Value = "'1.015"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").Value = "'309.59"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("D7").Value = "'0.4770"
$ws.Range("E7").Value = "  +2.40%  "
$ws.Range("D8").Value = "'0.3690"
$ws.Range("E8").Value = "  +1.18%  "
$ws.Range("D9").Value = "'0.07253"
$ws.Range("E9").Value = "  +1.83%  "
$ws.Range("D10").Value = "'0.9331"
$ws.Range("E10").Value = "  +2.28%  "
$ws.Range("D11").Value = "'19.89"
$ws.Range("E11").Value = "  +2.03%  "
$ws.Range("D12").Value = "'0.07799"
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("D13").Value = "'1.806.98"
$ws.Range("E13").Value = "  -3.38%  "
$ws.Range("D14").Value = "'5.396"
$ws.Range("E14").Value = "  +2.36%  "
$ws.Range("D15").Value = "'6.490"
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").Value = "'89.40"
$ws.Range("E16").Value = "  +2.01%  "
$ws.Range("D17").Value = "'1.018"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").Value = "'0.000008707"
$ws.Range("E18").Value = "  +1.55%  "
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("D20").Value = "'27.140.37"
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("D21").Value = "'14.62"
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("D22").Value = "'5.062"
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("D23").Value = "'10.66"
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("D25").Value = "'153.18"
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("D26").Value = "'18.37"
$ws.Range("E26").Value = "  +0.97%  "
$ws.Range("D27").Value = "'1.989"
$ws.Range("E27").Value = "  -1.51%  "
$ws.Range("D28").Value = "'114.77"
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("D29").Value = "'4.924"
$ws.Range("E29").Value = "  +1.02%  "
$ws.Range("D30").Value = "'0.08879"
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("D31").Value = "'3.304"
$ws.Range("E31").Value = "  +3.18%  "
$ws.Range("D32").Value = "'1.184"
$ws.Range("E32").Value = "  +1.34%  "
$ws.Range("D33").Value = "'4.527"
$ws.Range("E33").Value = "  +1.69%  "
$ws.Range("D34").Value = "'0.7383"
$ws.Range("E34").Value = "  -0.55%  "
$ws.Range("D35").Value = "'2.688"
$ws.Range("E35").Value = "  -3.44%  "
$ws.Range("D36").Value = "'1.116"
$ws.Range("E36").Value = "  +3.30%  "
$ws.Range("D37").Value = "'0.01978"
$ws.Range("D38").Value = "'0.05275"
$ws.Range("E38").Value = "  +2.20%  "
$ws.Range("D39").Value = "'2.972"
$ws.Range("E39").Value = "  +0.37%  "
$ws.Range("D40").Value = "'0.5285"
$ws.Range("E40").Value = "  +2.14%  "
$ws.Range("D41").Value = "'7.058"
$ws.Range("E41").Value = "  +2.37%  "
$ws.Range("D42").Value = "'0.1526"
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("D43").Value = "'8.310"
$ws.Range("E43").Value = "  +2.41%  "
$ws.Range("D44").Value = "'10.60"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("D45").Value = "'0.4748"
$ws.Range("E45").Value = "  +1.52%  "
$ws.Range("D46").Value = "'1.016"
$ws.Range("E46").Value = "  +0.76%  "
$ws.Range("D47").Value = "'102.17"
$ws.Range("E47").Value = "  +1.82%  "
$ws.Range("D48").Value = "'1.620"
$ws.Range("E48").Value = "  +1.34%  "
$ws.Range("D49").Value = "'66.04"
$ws.Range("E49").Value = "  +2.08%  "
$ws.Range("D50").Value = "'0.06063"
$ws.Range("E50").Value = "  +0.58%  "
$ws.Range("D51").Value = "'0.8943"
$ws.Range("E51").Value = "  +1.82%  "
